$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust column B width (closest attainable width to 14.453125)
$ws.Columns.Item(2).ColumnWidth = 13.6

# Update data values (rows 2-9, columns A-E) to the new recomputed figures
$values = @{
    2 = @(-0.48954910496300785, -0.4946532411966314,  2.1169021775963999,  0.24390213106938413, -4.1168784502176825)
    3 = @(-0.55158687719685295, -0.545685321865641,    2.366408849518074,   0.30030870120921832, -4.675491331345988)
    4 = @( 0.014567646015852356, 0.021094586376653356,-0.046004255836585406,0.035781695249631686, 0.039012422091630805)
    5 = @( 0.014676174674985778, 0.021708004314760592,-0.04715666409993139, 0.036465069727266666, 0.039101489229801258)
    6 = @( 1.6569998302275373,   1.5461160198025354,   0.75804350810178989,-6.5662890747631133,   20.206454679137028)
    7 = @( 1.8344092236928842,   1.7433671561146917,   0.85643559192981611,-7.3624174638324353,   22.712249239317536)
    8 = @(-0.098682227424718591, 0.0026991303881585037,-0.018205401421277162,0.041109985583532849, 0.069158423460431911)
    9 = @(-0.09964626092574011,  0.0032214427486452488,-0.019360154319327815,0.040913117377496512, 0.069510524450699732)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    for ($col = 1; $col -le 5; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
    }
}
